$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Ost/Nord (easting/northing) coordinates in row 2
$ws.Range("Q2").Value = 471365
$ws.Range("R2").Value = 6298471

# Remove the now-unused start/end time cells for row 2
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
